$d = $word.ActiveDocument

# Update the date line at the top of the document.
[void]$d.Content.Find.Execute("2023-09-28 Thursday", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "2023-09-29 Friday", 2)

# Update the division problems in the table. Cells are addressed directly by
# (row, column) rather than by text search because several of the values are
# duplicated and/or reused as both a source and a target string elsewhere in
# the table, which would make a global find/replace unsafe.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "39÷8="
$t.Cell(1, 2).Range.Text  = "80÷9="
$t.Cell(1, 3).Range.Text  = "67÷3="
$t.Cell(1, 4).Range.Text  = "95÷2="
$t.Cell(1, 5).Range.Text  = "13÷7="

$t.Cell(5, 1).Range.Text  = "56÷7="
$t.Cell(5, 2).Range.Text  = "58÷6="
$t.Cell(5, 3).Range.Text  = "10÷4="
$t.Cell(5, 4).Range.Text  = "15÷9="
$t.Cell(5, 5).Range.Text  = "95÷6="

$t.Cell(9, 1).Range.Text  = "39÷5="
$t.Cell(9, 2).Range.Text  = "40÷4="
$t.Cell(9, 3).Range.Text  = "46÷6="
$t.Cell(9, 4).Range.Text  = "55÷6="
$t.Cell(9, 5).Range.Text  = "58÷4="

$t.Cell(13, 1).Range.Text = "69÷9="
$t.Cell(13, 2).Range.Text = "77÷6="
$t.Cell(13, 3).Range.Text = "52÷4="
$t.Cell(13, 4).Range.Text = "53÷5="
$t.Cell(13, 5).Range.Text = "14÷4="

$t.Cell(17, 1).Range.Text = "59÷3="
$t.Cell(17, 2).Range.Text = "33÷2="
$t.Cell(17, 3).Range.Text = "44÷8="
$t.Cell(17, 4).Range.Text = "16÷8="
$t.Cell(17, 5).Range.Text = "47÷9="
